# Update "想去人数" (want-to-go count) figures for the 合肥-漫展信息 workbook
# to match the freshly generated data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 86
$wsExhibit.Range("F5").Value = 4788
$wsExhibit.Range("F10").Value = 215

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 24

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 86
$wsAll.Range("F5").Value = 4788
$wsAll.Range("F10").Value = 24
$wsAll.Range("F11").Value = 215
